# HRMS SQL schema sheet update
# Refines the raw field/table names into fully annotated SQL-style
# column definitions and fixes several naming inconsistencies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Users / Employees / Candidates / Employers (header block) ----
$ws.Range("A2").Value = "id(int,identity)"
$ws.Range("C2").Value = "id(int, primary key, foreign key)"
$ws.Range("E2").Value = "id(int, primary key, foreign key)"
$ws.Range("G2").Value = "id"

$ws.Range("A3").Value = "email(string, unique, not null)"
$ws.Range("C3").Value = "first_name(string, not null)"
$ws.Range("E3").Value = "first_name(string, not null)"
$ws.Range("G3").Value = "company_name"

$ws.Range("A4").Value = "password(string, not null)"
$ws.Range("C4").Value = "last_name(string, not null)"
$ws.Range("E4").Value = "last_name(string, not null)"
$ws.Range("G4").Value = "web_address"

# Employees column (C) no longer has a 3rd/4th row entry
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()

$ws.Range("E5").Value = "identity_number"
$ws.Range("G5").Value = "phone_number"

$ws.Range("E6").Value = "birth_year"

# ---- VerificationCodes / VerificationCodeCandidates / VerificationCodeEmloyeers ----
$ws.Range("A12").Value = "verification_codes"
$ws.Range("C12").Value = "erification_code_candidates"
$ws.Range("E12").Value = "Verification_code_emloyeers"

$ws.Range("C13").Value = "id"
$ws.Range("E13").Value = "id"

$ws.Range("C14").Value = "candidate_id"
$ws.Range("E14").Value = "employer_id"

# ---- EmployeeConfirms / (new) employee_confirm_employers / JobTitles ----
$ws.Range("A22").Value = "employee_confirms"
$ws.Range("C22").Value = "employee_confirm_employers"
$ws.Range("E22").Value = "job_titels"

$ws.Range("C23").Value = "id"

$ws.Range("A24").Value = "employee_id"
$ws.Range("C24").Value = "employer_id"
$ws.Range("E24").Value = "title"

$ws.Range("A25").Value = "is_confirmed"
$ws.Range("A26").Value = "confirm_date"

# Restore the cursor position recorded in the saved file
$ws.Range("E11").Select()
